# SCE_mean.xlsx - "Add files via upload / Update 05/2025"
# Refresh a handful of previously-computed mean values (tiny floating-point
# precision differences from recomputation), fill in two previously-blank
# data rows (2024-06-01 / 2024-07-01), and append two new trailing rows
# for 2025-03-01 and 2025-04-01 (dates only, values not yet available).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Tiny re-computation precision refreshes in column B (and one in C) ---
$ws.Range("B7").Value   = 5.3768105768561369
$ws.Range("B48").Value  = 2.8627888042811933
$ws.Range("B52").Value  = 4.4007252819400584
$ws.Range("B63").Value  = 4.3152441953006555
$ws.Range("C65").Value  = 3.9000000953674316
$ws.Range("B89").Value  = 4.8222522532140601
$ws.Range("B101").Value = 5.5214701284176035
$ws.Range("B104").Value = 6.5096316414125805
$ws.Range("B105").Value = 5.8315007289581144
$ws.Range("B106").Value = 7.003945672164912
$ws.Range("B121").Value = 6.6232865476542875

# --- Fill in the two rows that previously had no mean/median values yet ---
$ws.Range("B139").Value = 3.0038180497508256
$ws.Range("C139").Value = 3.5999999046325684

$ws.Range("B140").Value = 4.7476770051488746
$ws.Range("C140").Value = 3

# --- Append two new trailing rows (dates only; values not yet published) ---
# Copy formatting from the last existing row (147) down into the two new
# rows so the new date cells get the date number format and the new B/C
# cells get the same numeric format as the rest of the column.
$ws.Range("A147:C147").Copy()
$ws.Range("A148:C149").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(148, 1).Value = 45717
$ws.Cells.Item(149, 1).Value = 45748
